$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add "Mean increase" / "Median increase" headers (bold, like the
# other section headers) and the increase formulas for each of the
# Low / Medium / High / All blocks.
# ---------------------------------------------------------------

# Low block (rows 2-31 data, summary in rows 2-16)
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

$ws.Range("D19").Formula = "= ((E3 / 114.202998) * 100) - 100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").Formula = "= ((E10 / 113.658804) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Medium block (rows 34-63 data, summary in rows 34-48)
$ws.Range("D50").Value = "Mean increase"
$ws.Range("D50").Font.Bold = $true
$ws.Range("F50").Value = "Median increase"
$ws.Range("F50").Font.Bold = $true

$ws.Range("D51").Formula = "= ((E35 / 114.202998) * 100) - 100"
$ws.Range("D51").ClearFormats()
$ws.Range("F51").Formula = "= ((E42 / 113.658804) * 100) - 100"
$ws.Range("F51").ClearFormats()

# High block (rows 66-95 data, summary in rows 66-80)
$ws.Range("D82").Value = "Mean increase"
$ws.Range("D82").Font.Bold = $true
$ws.Range("F82").Value = "Median increase"
$ws.Range("F82").Font.Bold = $true

$ws.Range("D83").Formula = "= ((E67 / 114.202998) * 100) - 100"
$ws.Range("D83").ClearFormats()
$ws.Range("F83").Formula = "= ((E74 / 113.658804) * 100) - 100"
$ws.Range("F83").ClearFormats()

# All block (rows 98-187 data)
$ws.Range("D113").Value = "Mean increase"
$ws.Range("D113").Font.Bold = $true
$ws.Range("F113").Value = "Median increase"
$ws.Range("F113").Font.Bold = $true

$ws.Range("D114").Formula = "= (D19 + D51 + D83) / 3"
$ws.Range("D114").ClearFormats()
$ws.Range("F114").Formula = "= (F19 + F51 + F83) / 3"
$ws.Range("F114").ClearFormats()
